# amova_combined_table.docx edit script
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1. widen the first grid column: 4777 dxa (238.85pt) -> 5045 dxa (252.25pt)
$t.Columns.Item(1).Width = 252.25

# 2. "Among urban/rural groups" row: p 0.375 -> 0.336
$t.Cell(2, 7).Range.Text = "0.336"

# 3. "Among populations within urban/rural groups" -> "Among sampling sites within urban/rural groups"
$t.Cell(3, 1).Range.Text = "Among sampling sites within urban/rural groups"
# ... and its p value 0.001 -> (blank)
$t.Cell(3, 7).Range.Text = ""

# 4. "Within populations" -> "Within sampling sites"
$t.Cell(4, 1).Range.Text = "Within sampling sites"
# ... and its p value 0.001 -> (blank)
$t.Cell(4, 7).Range.Text = ""
# ... row height grows from 612 (30.6pt) to 614 (30.7pt)
$t.Rows.Item(4).Height = 30.7

# 5. insert two new rows before the old "Total" row (row 5). Row handles in
#    this object model are position-anchored, so re-resolve rows by fresh
#    index after every insert instead of reusing old variables.
$new1 = $t.Rows.Add($t.Rows.Item(5))
$new2 = $t.Rows.Add($t.Rows.Item(6))

# row 5 (new1): duplicate of the current "Total" figures
$t.Rows.Item(5).Height = 30.7
$t.Cell(5, 1).Range.Text = "Total"
$t.Cell(5, 2).Range.Text = "255"
$t.Cell(5, 3).Range.Text = "17,673.721"
$t.Cell(5, 4).Range.Text = "69.309"
$t.Cell(5, 5).Range.Text = "69.410"
$t.Cell(5, 6).Range.Text = "100%"
$t.Cell(5, 7).Range.Text = ""

# row 6 (new2): new "Variation within sampling sites" row
$t.Rows.Item(6).Height = 30.7
$t.Cell(6, 1).Range.Text = "Variation within sampling sites"
$t.Cell(6, 2).Range.Text = ""
$t.Cell(6, 3).Range.Text = ""
$t.Cell(6, 4).Range.Text = ""
$t.Cell(6, 5).Range.Text = ""
$t.Cell(6, 6).Range.Text = ""
$t.Cell(6, 7).Range.Text = "0.001"

# row 7 (the original "Total" row, shifted down): becomes
# "Variation between sampling sites" - figures cleared, p-value populated
$t.Rows.Item(7).Height = 30.7
$t.Cell(7, 1).Range.Text = "Variation between sampling sites"
$t.Cell(7, 2).Range.Text = ""
$t.Cell(7, 3).Range.Text = ""
$t.Cell(7, 4).Range.Text = ""
$t.Cell(7, 5).Range.Text = ""
$t.Cell(7, 6).Range.Text = ""
$t.Cell(7, 7).Range.Text = "0.001"

Write-Output "done"
